# Update the "timestamp" column (Z) on the Log_Muestras sheet with
# fresh run timestamps, as produced by a new execution of the logging
# script (dataset Us Crime agregado).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log_Muestras")

$ws.Range("Z2:Z5").Value = "2025-11-13T06:52:55.953880"
$ws.Range("Z6:Z9").Value = "2025-11-13T06:52:55.954877"
$ws.Range("Z10:Z16").Value = "2025-11-13T06:52:55.955880"
$ws.Range("Z17:Z26").Value = "2025-11-13T06:52:55.956877"
$ws.Range("Z27:Z35").Value = "2025-11-13T06:52:55.957873"
$ws.Range("Z36:Z45").Value = "2025-11-13T06:52:55.958725"
$ws.Range("Z46").Value = "2025-11-13T06:52:56.358123"
$ws.Range("Z47:Z69").Value = "2025-11-13T06:52:56.358661"
$ws.Range("Z70").Value = "2025-11-13T06:52:56.368977"
$ws.Range("Z71:Z74").Value = "2025-11-13T06:52:56.369489"
$ws.Range("Z75:Z76").Value = "2025-11-13T06:52:56.579351"
$ws.Range("Z77:Z86").Value = "2025-11-13T06:52:56.580352"
$ws.Range("Z87:Z95").Value = "2025-11-13T06:52:56.581352"
$ws.Range("Z96:Z102").Value = "2025-11-13T06:52:56.582351"
